$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 40
$ws.Range("B2").Value = 48
$ws.Range("B3").Value = 85
$ws.Range("B4").Value = 95
$ws.Range("B5").Value = 112
$ws.Range("B6").Value = 133
$ws.Range("B7").Value = 151
$ws.Range("B8").Value = 176
